# "H4MPT pathway in and precursors partway"
#
# 1. Rename Sheet2 -> H4MPT and populate it with the new H4MPT synthesis
#    pathway table (rows 1-8) plus a partial precursor-compound row (18-19).
# 2. Tidy up the selections on both sheets and make H4MPT the active tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Coenzyme B")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- rename the (until now empty) second sheet ---------------------------
$ws2.Name = "H4MPT"

# --- header row, reusing the same column headers as the Coenzyme B sheet -
$ws2.Range("A1").Value = "Reactions (normal form)"
$ws2.Range("B1").Value = "Rxn ID (Kbase)"
$ws2.Range("C1").Value = "Gene Name"
$ws2.Range("D1").Value = "Gene locus"

# --- the new strings are written in the exact order the author first ----
# --- typed them so the shared-string table grows the same way -----------

# the partial precursor-compound line (row 19) was started first
$ws2.Range("D19").Value = "C9H11N5O6P"
$ws2.Range("B19").Value = "cpd15850"

# column C (gene names) for rows 2,3,5,6,7,8
$ws2.Range("C2").Value = "mptA"
$ws2.Range("C3").Value = "mptB"
$ws2.Range("C5").Value = "mptD"
$ws2.Range("C6").Value = "mptE"
$ws2.Range("C7").Value = "mptG"
$ws2.Range("C8").Value = "mptH"

# column D (gene locus) for rows 2,3
$ws2.Range("D2").Value = "MMP0034"
$ws2.Range("D3").Value = "MMP0230"

# column B (rxn IDs) for rows 3-8
$ws2.Range("B3").Value = "rxn10490"
$ws2.Range("B4").Value = "rxn03168"
$ws2.Range("B5").Value = "rxn02504"
$ws2.Range("B6").Value = "rxn02503"
$ws2.Range("B7").Value = "rxn10446"
$ws2.Range("B8").Value = "rn10491"

# column D (gene locus) for rows 5-7
$ws2.Range("D5").Value = "MMP0243"
$ws2.Range("D6").Value = "MMP0579"
$ws2.Range("D7").Value = "MMP0279 (maybe)"

# column A (reaction text) for rows 2-5, then 7, then 6
$ws2.Range("A2").Value = "GTP + 2 H2O -> Formate + Ppi + 7,8-dihydronepterin 2' :3'-cyclicphosphate"
$ws2.Range("A3").Value = "7,8-dihydronepterin 2' :3'-cyclicphosphate + H2O -> Dihydroneopterin phosphate + H+"
$ws2.Range("A4").Value = " Dihydroneopterin phosphate + H2O <=> Dihydroneopterin + H+ + Ppi"
$ws2.Range("A5").Value = "Dihydroneopterin  -> 6-hydroxymethyl-7,8-dihydropterin  + Glycolaldehyde"
$ws2.Range("A7").Value = "4-aminobenzoate + PRPP -> beta-RFA-P"
$ws2.Range("A6").Value = "6-hydroxymethyl-7,8-dihydropterin + ATP -> 6-hydroxymethyl-7,8-dihydropterin diphosphate + Ppi"

# --- the compound-info sub header + partial row below the pathway table -
$ws2.Range("A18").Value = "Compound Name in Kbase (name above, if applicable)"
$ws2.Range("B18").Value = "Kbase ID"
$ws2.Range("C18").Value = "KEGG ID"
$ws2.Range("D18").Value = "Formula"
$ws2.Range("E18").Value = "Charge"
$ws2.Range("E19").Value = -1

# --- match the Courier-New "data" style used elsewhere in the workbook --
$ws1.Range("A25").Copy()
$ws2.Range("B19").PasteSpecial(-4122)
$ws2.Range("D19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- size the columns to fit their new content ---------------------------
$ws2.Columns.Item(1).AutoFit()
$ws2.Columns.Item(2).AutoFit()
$ws2.Columns.Item(3).AutoFit()
$ws2.Columns.Item(4).AutoFit()
$ws2.Columns.Item(5).AutoFit()

$ws2.PageSetup.Orientation = 1

# --- fix up the selections: Coenzyme B no longer needs its old D7 -------
# --- selection, and H4MPT becomes the active, selected sheet ------------
$ws1.Activate()
$ws1.Range("A1:XFD1").Select()

$ws2.Activate()
$ws2.Range("A13").Select()
